$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 1568
$ws.Range("E4").Value = 6597
$ws.Range("D5").Value = "659 (42.0)"
$ws.Range("E5").Value = "2526 (38.3)"
$ws.Range("D6").Value = "377 (24.0)"
$ws.Range("E6").Value = "1509 (22.9)"
$ws.Range("D7").Value = "318 (20.3)"
$ws.Range("E7").Value = "1456 (22.1)"
$ws.Range("D8").Value = "214 (13.6)"
$ws.Range("E8").Value = "1106 (16.8)"
$ws.Range("D9").Value = "731 (46.6)"
$ws.Range("E9").Value = "2786 (42.2)"
$ws.Range("D10").Value = "590 (37.6)"
$ws.Range("E10").Value = "3258 (49.4)"
$ws.Range("D11").Value = "222 (14.2)"
$ws.Range("E11").Value = "376 (5.7)"
$ws.Range("D12").Value = "756 (48.2)"
$ws.Range("E12").Value = "2963 (44.9)"
$ws.Range("D13").Value = "495 (31.6)"
$ws.Range("E13").Value = "336 (5.1)"
$ws.Range("D14").Value = "104 (6.6)"
$ws.Range("E14").Value = "753 (11.4)"
$ws.Range("D15").Value = "539 (34.4)"
$ws.Range("E15").Value = "2760 (41.8)"
$ws.Range("D16").Value = "365 (23.3)"
$ws.Range("E16").Value = "1508 (22.9)"
$ws.Range("D17").Value = "505 (32.2)"
$ws.Range("E17").Value = "2122 (32.2)"
$ws.Range("D18").Value = "1041 (66.4)"
$ws.Range("E18").Value = "4566 (69.2)"
$ws.Range("D19").Value = "245 (15.6)"
$ws.Range("E19").Value = "680 (10.3)"
$ws.Range("D20").Value = "819 (52.2)"
$ws.Range("E20").Value = "3823 (58.0)"
$ws.Range("D21").Value = "1131 (72.1)"
$ws.Range("E21").Value = "4913 (74.5)"
$ws.Range("D22").Value = "324 (20.7)"
$ws.Range("E22").Value = "976 (14.8)"
$ws.Range("D23").Value = "922 (58.8)"
$ws.Range("E23").Value = "4241 (64.3)"
$ws.Range("D24").Value = "1074 (68.5)"
$ws.Range("E24").Value = "4362 (66.1)"
$ws.Range("D25").Value = "199 (12.7)"
$ws.Range("E25").Value = "962 (14.6)"
$ws.Range("D26").Value = "1550 (98.9)"
$ws.Range("E26").Value = "6545 (99.2)"
$ws.Range("D27").Value = "1119 (71.4)"
$ws.Range("E27").Value = "4255 (64.5)"
$ws.Range("D28").Value = "623 (39.7)"
$ws.Range("E28").Value = "2549 (38.6)"
$ws.Range("D29").Value = "382 (24.4)"
$ws.Range("E29").Value = "1690 (25.6)"
$ws.Range("D30").Value = "21 (1.3)"
$ws.Range("E30").Value = "119 (1.8)"
$ws.Range("D31").Value = "463 (29.5)"
$ws.Range("E31").Value = "2432 (36.9)"
$ws.Range("D33").Value = "10 (0.6)"
$ws.Range("E33").Value = "43 (0.7)"
$ws.Range("D34").Value = "121 (7.7)"
$ws.Range("E34").Value = "312 (4.7)"
$ws.Range("D35").Value = "44 (2.8)"
$ws.Range("E35").Value = "124 (1.9)"
$ws.Range("D36").Value = "185 (11.8)"
$ws.Range("E36").Value = "338 (5.1)"
$ws.Range("D37").Value = "1207 (77.0)"
$ws.Range("E37").Value = "5780 (87.6)"
$ws.Range("D38").Value = "49 (3.1)"
$ws.Range("E38").Value = "179 (2.7)"
$ws.Range("D39").Value = "634 (40.4)"
$ws.Range("E39").Value = "1828 (27.7)"
$ws.Range("D40").Value = "885 (56.4)"
$ws.Range("E40").Value = "4590 (69.6)"
$ws.Range("D41").Value = "73 (4.7)"
$ws.Range("E41").Value = "315 (4.8)"
$ws.Range("D42").Value = "90 (5.7)"
$ws.Range("E42").Value = "358 (5.4)"
$ws.Range("E43").Value = "42 (0.6)"
$ws.Range("E44").Value = "13 (0.2)"
$ws.Range("D46").Value = "64 [51,75]"
$ws.Range("E46").Value = "68 [57,78]"
$ws.Range("C47").Value = 6292
$ws.Range("D47").Value = "8.00 [5.50,12.58]"
$ws.Range("E47").Value = "8.25 [5.70,13.00]"
$ws.Range("C48").Value = 1873
$ws.Range("D48").Value = "7.54 [5.13,12.58]"
$ws.Range("E48").Value = "7.17 [5.13,11.79]"
$ws.Range("C49").Value = 6292
$ws.Range("C50").Value = 1873
$ws.Range("D50").Value = "17.00 [11.00,26.00]"
$ws.Range("E50").Value = "15.00 [10.00,24.00]"
$ws.Range("D51").Value = "6 [4,9]"
$ws.Range("C53").Value = 2409
$ws.Range("C54").Value = 16
$ws.Range("C55").Value = 2934
$ws.Range("C56").Value = 21
$ws.Range("C57").Value = 26
$ws.Range("C59").Value = 3043
$ws.Range("D59").Value = "915 [289,1601]"
$ws.Range("E59").Value = "877 [285,1695]"
$ws.Range("D60").Value = "4722 [2123,9457]"
$ws.Range("E60").Value = "4403 [1879,8700]"
$ws.Range("D61").Value = "553.1 [296.4,947.0]"
$ws.Range("E61").Value = "516.5 [265.4,884.1]"
$ws.Range("C62").Value = 4359
$ws.Range("D62").Value = "50 [40,60]"
$ws.Range("C63").Value = 2121
$ws.Range("D63").Value = "66.0 [27.0,134.0]"
$ws.Range("E63").Value = "56.0 [23.0,116.0]"
$ws.Range("C64").Value = 2121
$ws.Range("D64").Value = "0.32 [0.15,0.52]"
$ws.Range("E64").Value = "0.27 [0.13,0.47]"
$ws.Range("C65").Value = 2121
$ws.Range("D65").Value = "3.0 [1.0,13.0]"
$ws.Range("E65").Value = "3.0 [1.0,14.0]"
$ws.Range("C66").Value = 6865
$ws.Range("D66").Value = "24.0 [5.0,69.5]"
$ws.Range("E66").Value = "36.0 [8.0,88.0]"
$ws.Range("C67").Value = 2952
$ws.Range("D67").Value = "4.0 [1.0,22.2]"
$ws.Range("E67").Value = "4.0 [1.0,20.0]"
$ws.Range("C68").Value = 2952
$ws.Range("D68").Value = "48.5 [15.0,104.0]"
$ws.Range("E68").Value = "49.0 [18.0,101.0]"
$ws.Range("C69").Value = 2952
$ws.Range("D69").Value = "0.22 [0.07,0.49]"
$ws.Range("E69").Value = "0.24 [0.09,0.47]"
$ws.Range("C70").Value = 26
$ws.Range("D70").Value = "19.7 [17.1,23.0]"
$ws.Range("E70").Value = "19.5 [17.1,22.5]"
$ws.Range("C71").Value = 21
$ws.Range("D71").Value = "77.6 [71.3,86.0]"
$ws.Range("E71").Value = "75.1 [69.4,82.2]"
$ws.Range("C72").Value = 460
$ws.Range("C73").Value = 24
$ws.Range("D73").Value = "97.8 [96.1,99.2]"
$ws.Range("E73").Value = "97.2 [95.7,98.6]"
$ws.Range("C74").Value = 21
$ws.Range("D74").Value = "88.5 [76.9,101.4]"
$ws.Range("E74").Value = "87.3 [76.2,99.6]"
$ws.Range("C75").Value = 2018
$ws.Range("D75").Value = "84.0 [66.0,118.0]"
$ws.Range("C76").Value = 2018
$ws.Range("D76").Value = "44.0 [37.0,54.0]"
$ws.Range("C77").Value = 1140
$ws.Range("C78").Value = 30
$ws.Range("D78").Value = "160.0 [124.0,224.5]"
$ws.Range("E78").Value = "152.0 [123.0,201.0]"
$ws.Range("C79").Value = 13
$ws.Range("C80").Value = 17
$ws.Range("C81").Value = 7958
$ws.Range("D81").Value = "18.1 [13.1,30.4]"
$ws.Range("E81").Value = "23.7 [13.6,37.7]"
$ws.Range("C82").Value = 1089
$ws.Range("C83").Value = 5303
$ws.Range("D83").Value = "231.0 [150.5,361.5]"
$ws.Range("E83").Value = "237.0 [160.5,363.5]"
$ws.Range("C84").Value = 458
